$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.274.12"
$ws.Range("E2").Value = "  +5.45%  "
$ws.Range("D3").Value = "1.916.99"
$ws.Range("E3").Value = "  +5.89%  "
$ws.Range("D5").Value = "'254.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.5155"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.74%  "
$ws.Range("D8").Value = "'46.04"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +6.78%  "
$ws.Range("D9").Value = "'0.2975"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.12%  "
$ws.Range("D10").Value = "'0.06842"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.96%  "
$ws.Range("D11").Value = "1.915.94"
$ws.Range("E11").Value = "  +5.90%  "
$ws.Range("D12").Value = "'17.51"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.33%  "
$ws.Range("D13").Value = "'0.07336"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.19%  "
$ws.Range("D14").Value = "'0.6899"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.05%  "
$ws.Range("D15").Value = "'87.62"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.89%  "
$ws.Range("D16").Value = "'4.914"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.11%  "
$ws.Range("D17").Value = "30.274.73"
$ws.Range("E17").Value = "  +5.57%  "
$ws.Range("D18").Value = "'0.000008045"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +9.28%  "
$ws.Range("D19").Value = "'0.9991"
$ws.Range("D19").ClearFormats()
$ws.Range("E20").Value = "  +6.54%  "
$ws.Range("D21").Value = "2.164.78"
$ws.Range("E21").Value = "  +6.29%  "
$ws.Range("D22").Value = "'0.9987"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'4.859"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +5.15%  "
$ws.Range("D24").Value = "'5.767"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +8.41%  "
$ws.Range("D25").Value = "'9.192"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("D28").Value = "'17.29"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.92%  "
$ws.Range("D29").Value = "'2.014"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +6.61%  "
$ws.Range("D30").Value = "'1.377"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").Value = "'4.275"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").Value = "'0.08853"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.70%  "
$ws.Range("D33").Value = "'4.041"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.10%  "
$ws.Range("D34").Value = "'0.05132"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("D35").Value = "'1.159"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.06%  "
$ws.Range("D36").Value = "'0.7196"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.60%  "
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("D38").Value = "'2.328"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.85%  "
$ws.Range("D39").Value = "'2.836"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.70%  "
$ws.Range("D40").Value = "'0.9772"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").Value = "'0.01699"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.62%  "
$ws.Range("D42").Value = "'6.134"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("D43").Value = "'0.4331"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.00%  "
$ws.Range("D44").Value = "'105.71"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.63%  "
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'7.704"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.43%  "
$ws.Range("D47").Value = "'0.1278"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.21%  "
$ws.Range("D48").Value = "'0.05732"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.36%  "
$ws.Range("D49").Value = "'8.561"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.61%  "
$ws.Range("D50").Value = "'33.41"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.29%  "
$ws.Range("E51").Value = "  +6.71%  "

# Row 26/27 swap (Monero <-> BitcoinCash)
$ws.Range("B26").Value = "BitcoinCash"
$ws.Range("C26").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D26").Value = "'140.01"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +24.63%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'146.47"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.45%  "
